$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text formatting (values like "208.30"
# or "7.10" must not be auto-converted to numbers, which would drop
# trailing zeros).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.645.89'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.19%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.598.15'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.08%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.73'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.29%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.44%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.41%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.59'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.21%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.27%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.822.19'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.11%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.586.64'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.49%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.63'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.609.66'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.26%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.37%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '208.30'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.62%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.10'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.28'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.25%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.94'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.67'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.34%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.19%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.27'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.38%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.96%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.30%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.41%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.27'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +19.30%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.278.52'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.88%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.48'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.71%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.59%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0167'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.00%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.823'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.15%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.16'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.09%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.44%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '62.64'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.733.78'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '89.60'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.53%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.15%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.13%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0513'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.63%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.17%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.45'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.13%  '
